$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: SARBM3D wasn't loading from the C drive on ECS lab computers, so
# the BM3DELBP row (row 8) had a few results stuck on red "desktop"
# placeholder markers instead of the actual computed numbers. Now that it
# loads correctly there, fill in the real values.
#
# J8's correct formatting (bold, highlighted) is the same highlight style
# H8/U8 currently carry (they were only flagged because their SARBM3D-derived
# numbers looked wrong) - grab that formatting before H8/U8 get normalised
# to the regular "complete" look below.
$ws.Range("H8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Replace the placeholder text values with the real computed results.
$ws.Range("H8").Value = 0.93906250000000002
$ws.Range("J8").Value = 0.1140625
$ws.Range("L8").Value = 0.72834821428571395

# H8, L8 and U8 are now normal completed results, so match the formatting
# already used by the other finished cells in this row.
$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("U8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the view where the user was last working.
$null = $ws.Range("L8").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
